# Add drop_geo ablation outputs
$p = $ppt.ActivePresentation

# --- Slide 6 "Model Sonuclari" -> Content Placeholder 2 ---
# Rename the full-data model labels, then append the two No-Geo ablation rows.
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange

$s6p1 = $tr6.Paragraphs(1, 1)
$s6p1.Text = ""
$s6p1 = $tr6.Paragraphs(1, 1)
$s6p1.Text = "Full LR -> MAE 94.83, RMSE 232.56, R2 -0.115"

$s6p2 = $tr6.Paragraphs(2, 1)
$s6p2.Text = ""
$s6p2 = $tr6.Paragraphs(2, 1)
$s6p2.Text = "Full RF -> MAE 42.15, RMSE 156.90, R2 0.492"

$s6p3 = $tr6.Paragraphs(3, 1)
$s6p3.Text = ""
$s6p3 = $tr6.Paragraphs(3, 1)
$s6p3.Text = "No-Geo LR -> MAE 74.62, RMSE 211.84, R2 0.074"

$s6p3 = $tr6.Paragraphs(3, 1)
$null = $s6p3.InsertAfter("`rNo-Geo RF -> MAE 25.98, RMSE 116.89, R2 0.718")

# --- Slide 8 "Sonuc ve Kisa Ozet" -> Content Placeholder 2 ---
# Shift the tail two bullets: a new No-Geo RF summary line replaces the
# "Karlilik iliskileri..." slot, and that line moves down to replace the
# old 15-second-summary line.
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange

$s8p3 = $tr8.Paragraphs(3, 1)
$s8p3.Text = ""
$s8p3 = $tr8.Paragraphs(3, 1)
$s8p3.Text = "No-Geo RF ile R2 0.718 seviyesine cikti"

$s8p4 = $tr8.Paragraphs(4, 1)
$s8p4.Text = ""
$s8p4 = $tr8.Paragraphs(4, 1)
$s8p4.Text = "Karlilik iliskileri dogrusal degil; indirim ve kategori etkisi kritik"

# --- Slide 9 "Sinirlamalar ve Ileri Isler" -> Content Placeholder 2 ---
# Insert a new bullet right after the first one.
$s9 = $p.Slides.Item(9)
$tr9 = $s9.Shapes.Item(2).TextFrame.TextRange
$s9p1 = $tr9.Paragraphs(1, 1)
$null = $s9p1.InsertAfter("`rdrop_geo testi geo kolonlarin gurultu olabilecegini gosterdi")

# --- Slide 10 "Soru - Cevap" -> Content Placeholder 2 ---
# Shorten "Random Forest" to "RF" in the second bullet.
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(2).TextFrame.TextRange
$s10p2 = $tr10.Paragraphs(2, 1)
$s10p2.Text = ""
$s10p2 = $tr10.Paragraphs(2, 1)
$s10p2.Text = "Neden RF? -> Dogrusal olmayan iliskileri yakaliyor"
